$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Miesiac 1")
$ws.Range("B2").Value = 0.2316763795256406
$ws.Range("C2").Value = 0.2830683067390554
$ws.Range("D2").Value = 0.2144229011851627
$ws.Range("B3").Value = 0.273196450901372
$ws.Range("C3").Value = 0.2788866191901092
$ws.Range("D3").Value = 0.2840676161439734
$ws.Range("B4").Value = 0.2186905015119305
$ws.Range("C4").Value = 0.2278386555315274
$ws.Range("D4").Value = 0.2324696936659122

$ws = $wb.Worksheets.Item("Miesiac 2")
$ws.Range("B2").Value = 0.2316763795256406
$ws.Range("C2").Value = 0.2830683067390554
$ws.Range("D2").Value = 0.2144229011851627
$ws.Range("B3").Value = 0.273196450901372
$ws.Range("C3").Value = 0.2788866191901092
$ws.Range("D3").Value = 0.2840676161439734
$ws.Range("B4").Value = 0.2186905015119305
$ws.Range("C4").Value = 0.2278386555315274
$ws.Range("D4").Value = 0.2324696936659122

$ws = $wb.Worksheets.Item("Miesiac 3")
$ws.Range("B2").Value = 0.7722545984188022
$ws.Range("C2").Value = 0.9435610224635179
$ws.Range("D2").Value = 0.7147430039505425
$ws.Range("B3").Value = 0.9106548363379067
$ws.Range("C3").Value = 0.9296220639670305
$ws.Range("D3").Value = 0.9468920538132446
$ws.Range("B4").Value = 0.7289683383731018
$ws.Range("C4").Value = 0.7594621851050913
$ws.Range("D4").Value = 0.7748989788863739

$ws = $wb.Worksheets.Item("Miesiac 4")
$ws.Range("B2").Value = 0.5405782188931615
$ws.Range("C2").Value = 0.6604927157244626
$ws.Range("D2").Value = 0.5003201027653797
$ws.Range("B3").Value = 0.6374583854365347
$ws.Range("C3").Value = 0.6507354447769214
$ws.Range("D3").Value = 0.6628244376692711
$ws.Range("B4").Value = 0.5102778368611712
$ws.Range("C4").Value = 0.5316235295735638
$ws.Range("D4").Value = 0.5424292852204617

$ws = $wb.Worksheets.Item("Miesiac 5")
$ws.Range("B2").Value = 0.5405782188931615
$ws.Range("C2").Value = 0.6604927157244626
$ws.Range("D2").Value = 0.5003201027653797
$ws.Range("B3").Value = 0.6374583854365347
$ws.Range("C3").Value = 0.6507354447769214
$ws.Range("D3").Value = 0.6628244376692711
$ws.Range("B4").Value = 0.5102778368611712
$ws.Range("C4").Value = 0.5316235295735638
$ws.Range("D4").Value = 0.5424292852204617

$ws = $wb.Worksheets.Item("Miesiac 6")
$ws.Range("B2").Value = 0.5405782188931615
$ws.Range("C2").Value = 0.6604927157244626
$ws.Range("D2").Value = 0.5003201027653797
$ws.Range("B3").Value = 0.6374583854365347
$ws.Range("C3").Value = 0.6507354447769214
$ws.Range("D3").Value = 0.6628244376692711
$ws.Range("B4").Value = 0.5102778368611712
$ws.Range("C4").Value = 0.5316235295735638
$ws.Range("D4").Value = 0.5424292852204617

$ws = $wb.Worksheets.Item("Miesiac 7")
$ws.Range("B2").Value = 0.7722545984188022
$ws.Range("C2").Value = 0.9435610224635179
$ws.Range("D2").Value = 0.7147430039505425
$ws.Range("B3").Value = 0.9106548363379067
$ws.Range("C3").Value = 0.9296220639670305
$ws.Range("D3").Value = 0.9468920538132446
$ws.Range("B4").Value = 0.7289683383731018
$ws.Range("C4").Value = 0.7594621851050913
$ws.Range("D4").Value = 0.7748989788863739

$ws = $wb.Worksheets.Item("Miesiac 8")
$ws.Range("B2").Value = 0.5405782188931615
$ws.Range("C2").Value = 0.6604927157244626
$ws.Range("D2").Value = 0.5003201027653797
$ws.Range("B3").Value = 0.6374583854365347
$ws.Range("C3").Value = 0.6507354447769214
$ws.Range("D3").Value = 0.6628244376692711
$ws.Range("B4").Value = 0.5102778368611712
$ws.Range("C4").Value = 0.5316235295735638
$ws.Range("D4").Value = 0.5424292852204617

$ws = $wb.Worksheets.Item("Miesiac 9")
$ws.Range("B2").Value = 0.5405782188931615
$ws.Range("C2").Value = 0.6604927157244626
$ws.Range("D2").Value = 0.5003201027653797
$ws.Range("B3").Value = 0.6374583854365347
$ws.Range("C3").Value = 0.6507354447769214
$ws.Range("D3").Value = 0.6628244376692711
$ws.Range("B4").Value = 0.5102778368611712
$ws.Range("C4").Value = 0.5316235295735638
$ws.Range("D4").Value = 0.5424292852204617

$ws = $wb.Worksheets.Item("Miesiac 10")
$ws.Range("B2").Value = 0.5405782188931615
$ws.Range("C2").Value = 0.6604927157244626
$ws.Range("D2").Value = 0.5003201027653797
$ws.Range("B3").Value = 0.6374583854365347
$ws.Range("C3").Value = 0.6507354447769214
$ws.Range("D3").Value = 0.6628244376692711
$ws.Range("B4").Value = 0.5102778368611712
$ws.Range("C4").Value = 0.5316235295735638
$ws.Range("D4").Value = 0.5424292852204617

$ws = $wb.Worksheets.Item("Miesiac 11")
$ws.Range("B2").Value = 0.7722545984188022
$ws.Range("C2").Value = 0.9435610224635179
$ws.Range("D2").Value = 0.7147430039505425
$ws.Range("B3").Value = 0.9106548363379067
$ws.Range("C3").Value = 0.9296220639670305
$ws.Range("D3").Value = 0.9468920538132446
$ws.Range("B4").Value = 0.7289683383731018
$ws.Range("C4").Value = 0.7594621851050913
$ws.Range("D4").Value = 0.7748989788863739

$ws = $wb.Worksheets.Item("Miesiac 12")
$ws.Range("B2").Value = 0.2316763795256406
$ws.Range("C2").Value = 0.2830683067390554
$ws.Range("D2").Value = 0.2144229011851627
$ws.Range("B3").Value = 0.273196450901372
$ws.Range("C3").Value = 0.2788866191901092
$ws.Range("D3").Value = 0.2840676161439734
$ws.Range("B4").Value = 0.2186905015119305
$ws.Range("C4").Value = 0.2278386555315274
$ws.Range("D4").Value = 0.2324696936659122
